$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 167.75
$ws.Range("I9").Value = 101
$ws.Range("J9").Value = 234.5
$ws.Range("K9").Value = 101
$ws.Range("L9").Value = 234.5
$ws.Range("M9").Value = 68
$ws.Range("N9").Value = -572.5

$ws.Range("H34").Value = 3499.6667
$ws.Range("I34").Value = 3499.6667
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3499.6667
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3296.6667
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 3499.6667
$ws.Range("I36").Value = 3499.6667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 3499.6667
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2784.6667
$ws.Range("N36").ClearContents()

$ws.Range("H51").Value = 7328.5713
$ws.Range("I51").Value = 8325
$ws.Range("J51").Value = 6000
$ws.Range("K51").Value = 8325
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = -7841
$ws.Range("N51").Value = -6968

$ws.Range("H117").Value = 67998.336
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 67998.336
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 67998.336
$ws.Range("N117").Value = -77176.336

$ws.Range("H137").Value = 1037.4445
$ws.Range("I137").Value = 973.6875
$ws.Range("J137").Value = 1547.5
$ws.Range("K137").Value = 2921.0625
$ws.Range("L137").Value = 4642.5
$ws.Range("M137").Value = -371.0625
$ws.Range("N137").Value = -9742.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6612.222
$ws.Range("I2").Value = 6500
$ws.Range("J2").Value = 6702
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 6702
$ws.Range("M2").Value = -6387
$ws.Range("N2").Value = -6928

$ws.Range("H32").Value = 15744.444
$ws.Range("I32").Value = 14029.77
$ws.Range("J32").Value = 20202.6
$ws.Range("K32").Value = 14029.77
$ws.Range("L32").Value = 20202.6
$ws.Range("M32").Value = -13742.77
$ws.Range("N32").Value = -20776.6

$ws.Range("H61").Value = 5618.225
$ws.Range("I61").Value = 5729.763
$ws.Range("J61").Value = 3499
$ws.Range("K61").Value = 5729.763
$ws.Range("L61").Value = 3499
$ws.Range("M61").Value = -5517.763
$ws.Range("N61").Value = -3923

$ws.Range("H74").Value = 2683.65
$ws.Range("I74").Value = 2035.421
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 2035.421
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -1161.421
$ws.Range("N74").Value = -16748

$ws.Range("H77").Value = 2683.65
$ws.Range("I77").Value = 2035.421
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 10177.105
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -5809.105
$ws.Range("N77").Value = -83736

$ws.Range("H102").Value = 2716.5173
$ws.Range("I102").Value = 1951.16
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 1951.16
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -329.1600000000001
$ws.Range("N102").Value = -10744

$ws.Range("H116").Value = 6612.222
$ws.Range("I116").Value = 6500
$ws.Range("J116").Value = 6702
$ws.Range("K116").Value = 6500
$ws.Range("L116").Value = 6702
$ws.Range("M116").Value = -4206
$ws.Range("N116").Value = -11290

$ws.Range("H132").Value = 1579.775
$ws.Range("I132").Value = 1620.8948
$ws.Range("J132").Value = 798.5
$ws.Range("K132").Value = 4862.6844
$ws.Range("L132").Value = 2395.5
$ws.Range("M132").Value = -2332.6844
$ws.Range("N132").Value = -7455.5

$ws.Range("H136").Value = 5618.225
$ws.Range("I136").Value = 5729.763
$ws.Range("J136").Value = 3499
$ws.Range("K136").Value = 17189.289
$ws.Range("L136").Value = 10497
$ws.Range("M136").Value = -14639.289
$ws.Range("N136").Value = -15597

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6612.222
$ws.Range("I3").Value = 6500
$ws.Range("J3").Value = 6702
$ws.Range("K3").Value = 6500
$ws.Range("L3").Value = 6702
$ws.Range("M3").Value = -6386
$ws.Range("N3").Value = -6930

$ws.Range("H105").Value = 3081.6
$ws.Range("I105").Value = 3081.6
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3081.6
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1334.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 31000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 31000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 31000
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -33290

$ws.Range("H60").Value = 15989.223
$ws.Range("I60").Value = 18750
$ws.Range("J60").Value = 15200.429
$ws.Range("K60").Value = 18750
$ws.Range("L60").Value = 15200.429
$ws.Range("M60").Value = -18239
$ws.Range("N60").Value = -16222.429

$ws.Range("H62").Value = 4225
$ws.Range("I62").Value = 1905
$ws.Range("J62").Value = 4998.3335
$ws.Range("K62").Value = 1905
$ws.Range("L62").Value = 4998.3335
$ws.Range("M62").Value = -1281
$ws.Range("N62").Value = -6246.3335

$ws.Range("H65").Value = 4225
$ws.Range("I65").Value = 1905
$ws.Range("J65").Value = 4998.3335
$ws.Range("K65").Value = 9525
$ws.Range("L65").Value = 24991.6675
$ws.Range("M65").Value = -6405
$ws.Range("N65").Value = -31231.6675

$ws.Range("H74").Value = 37290
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 37290
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 37290
$ws.Range("N74").Value = -39038

$ws.Range("H77").Value = 37290
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 37290
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 111870
$ws.Range("N77").Value = -120606

$ws.Range("H95").Value = 26755.572
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 26755.572
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 26755.572
$ws.Range("N95").Value = -32247.572

$ws.Range("H141").Value = 317413.84
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 317413.84
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 317413.84
$ws.Range("N141").Value = -327773.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 224.66667
$ws.Range("I12").Value = 19.666666
$ws.Range("J12").Value = 361.33334
$ws.Range("K12").Value = 58.999998
$ws.Range("L12").Value = 1084.00002
$ws.Range("M12").Value = 114.000002
$ws.Range("N12").Value = -1430.00002

$ws.Range("H61").Value = 179.6
$ws.Range("I61").Value = 46.5
$ws.Range("J61").Value = 268.33334
$ws.Range("K61").Value = 139.5
$ws.Range("L61").Value = 805.0000200000001
$ws.Range("M61").Value = 75.5
$ws.Range("N61").Value = -1235.00002

$ws.Range("H87").Value = 8050.1
$ws.Range("I87").Value = 6166.778
$ws.Range("J87").Value = 25000
$ws.Range("K87").Value = 18500.334
$ws.Range("L87").Value = 75000
$ws.Range("M87").Value = -17252.334
$ws.Range("N87").Value = -77496

$ws.Range("H90").Value = 8050.1
$ws.Range("I90").Value = 6166.778
$ws.Range("J90").Value = 25000
$ws.Range("K90").Value = 55501.002
$ws.Range("L90").Value = 225000
$ws.Range("M90").Value = -49261.002
$ws.Range("N90").Value = -237480

$ws.Range("H121").Value = 126162.75
$ws.Range("I121").Value = 523
$ws.Range("J121").Value = 201546.6
$ws.Range("K121").Value = 1569
$ws.Range("L121").Value = 604639.8
$ws.Range("M121").Value = -259
$ws.Range("N121").Value = -607259.8

$ws.Range("H131").Value = 15153870
$ws.Range("I131").Value = 41667744
$ws.Range("J131").Value = 3083.5715
$ws.Range("K131").Value = 125003232
$ws.Range("L131").Value = 9250.7145
$ws.Range("M131").Value = -124998192
$ws.Range("N131").Value = -19330.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3938.111
$ws.Range("I14").Value = 1949.3334
$ws.Range("J14").Value = 4932.5
$ws.Range("K14").Value = 1949.3334
$ws.Range("L14").Value = 4932.5
$ws.Range("M14").Value = -1781.3334
$ws.Range("N14").Value = -5268.5

$ws.Range("H31").Value = 15000
$ws.Range("I31").Value = 15000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -14708

$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 15000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -14723

$ws.Range("H80").Value = 3150.7222
$ws.Range("I80").Value = 3084.2856
$ws.Range("J80").Value = 3193
$ws.Range("K80").Value = 3084.2856
$ws.Range("L80").Value = 3193
$ws.Range("M80").Value = -2086.2856
$ws.Range("N80").Value = -5189

$ws.Range("H83").Value = 3150.7222
$ws.Range("I83").Value = 3084.2856
$ws.Range("J83").Value = 3193
$ws.Range("K83").Value = 15421.428
$ws.Range("L83").Value = 15965
$ws.Range("M83").Value = -10429.428
$ws.Range("N83").Value = -25949

$ws.Range("H132").Value = 5477.364
$ws.Range("I132").Value = 5176.5557
$ws.Range("J132").Value = 6831
$ws.Range("K132").Value = 15529.6671
$ws.Range("L132").Value = 20493
$ws.Range("M132").Value = -12999.6671
$ws.Range("N132").Value = -25553

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4928.5557
$ws.Range("I7").Value = 4711.375
$ws.Range("J7").Value = 6666
$ws.Range("K7").Value = 4711.375
$ws.Range("L7").Value = 6666
$ws.Range("M7").Value = -4599.375
$ws.Range("N7").Value = -6890

$ws.Range("H126").Value = 4928.5557
$ws.Range("I126").Value = 4711.375
$ws.Range("J126").Value = 6666
$ws.Range("K126").Value = 14134.125
$ws.Range("L126").Value = 19998
$ws.Range("M126").Value = -11664.125
$ws.Range("N126").Value = -24938

$ws.Range("H136").Value = 50006784
$ws.Range("I136").Value = 31256736
$ws.Range("J136").Value = 125006984
$ws.Range("K136").Value = 93770208
$ws.Range("L136").Value = 375020952
$ws.Range("M136").Value = -93767658
$ws.Range("N136").Value = -375026052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 47877.555
$ws.Range("I54").Value = 22349.666
$ws.Range("J54").Value = 98933.336
$ws.Range("K54").Value = 22349.666
$ws.Range("L54").Value = 98933.336
$ws.Range("M54").Value = -21829.666
$ws.Range("N54").Value = -99973.336

$ws.Range("H136").Value = 4815.6763
$ws.Range("I136").Value = 4508.1377
$ws.Range("J136").Value = 6599.4
$ws.Range("K136").Value = 13524.4131
$ws.Range("L136").Value = 19798.2
$ws.Range("M136").Value = -10974.4131
$ws.Range("N136").Value = -24898.2
